# "case with 380 kV done" - update loading_percent results (rows 2-25, cols C-I & K).
# Column J remains 0 and is not touched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    @{ Row = 2; Values = @(3.213347429774682, 4.215992714421155, 16.51423181817644, 23.24235537693069, 27.43169869525295, 13.57650136409413, 19.31287144400861, 15.9273869945082) },
    @{ Row = 3; Values = @(3.132491027746396, 4.218149613489993, 15.57347455010789, 23.23689536000734, 27.35017530394033, 13.64743728563969, 19.3596527320048, 15.17854272015819) },
    @{ Row = 4; Values = @(3.081461755795277, 4.220022782119442, 14.97132330687864, 23.24580141741454, 27.319776176709, 13.6955607829372, 19.39733396430611, 14.70038695519085) },
    @{ Row = 5; Values = @(3.060340063165704, 4.220924455133294, 14.72004496828054, 23.25249209753605, 27.3123049493227, 13.71631219989344, 19.41492119192929, 14.50115389295924) },
    @{ Row = 6; Values = @(3.056813721355881, 4.221082542542677, 14.6779731669933, 23.25378728226805, 27.31136031890212, 13.71982661504088, 19.41797570274231, 14.46781445208795) },
    @{ Row = 7; Values = @(3.081178195962271, 4.220034381865823, 14.96795795980866, 23.24587928599649, 27.31965555325843, 13.69583603598905, 19.39756214374973, 14.69771742609175) },
    @{ Row = 8; Values = @(3.185767327438537, 4.216622703940785, 16.19509433423082, 23.23791708468208, 27.39949764183821, 13.60000700555339, 19.32713063970326, 15.67312753282774) },
    @{ Row = 9; Values = @(3.379045795360792, 4.214273299968677, 18.49171446928794, 23.32030272264496, 27.71259653514988, 13.44869670618302, 19.26092822817455, 17.4315470655445) },
    @{ Row = 10; Values = @(3.51281231663974, 4.215176406947814, 20.14526150089166, 23.4412642152516, 28.03795381007397, 13.36033852614759, 19.25715762090507, 18.61987276935374) },
    @{ Row = 11; Values = @(3.571689271291115, 4.216154635238645, 20.85559858042224, 23.50946446338764, 28.20638831633754, 13.32519909708403, 19.26536446991102, 19.13648997603784) },
    @{ Row = 12; Values = @(3.593686662351963, 4.216606339128972, 21.11860289578829, 23.53718119112255, 28.2730669775783, 13.31262810742824, 19.26991046241186, 19.32857757124554) },
    @{ Row = 13; Values = @(3.588962597008256, 4.216505447197892, 21.06222582204984, 23.53112791419763, 28.25857850121396, 13.31530264494207, 19.26886727689233, 19.28736708760336) },
    @{ Row = 14; Values = @(3.573505059291872, 4.216190170256352, 20.87735584545545, 23.51170687927954, 28.21181624207031, 13.32415008365941, 19.26570958935127, 19.15236467244832) },
    @{ Row = 15; Values = @(3.563997660910478, 4.216007628260096, 20.76333929584255, 23.50005697898374, 28.18354876063264, 13.32966543982447, 19.26396301743375, 19.06920760774544) },
    @{ Row = 16; Values = @(3.508923603969015, 4.215123867471225, 20.09799984871979, 23.43707206380253, 28.0273539996483, 13.36273743257144, 19.25682197133859, 18.58561965091852) },
    @{ Row = 17; Values = @(3.474621333883633, 4.21472680463684, 19.67913740266638, 23.40180600914114, 27.93673977418099, 13.38432694156394, 19.25499130523244, 18.28274165064099) },
    @{ Row = 18; Values = @(3.45470652205318, 4.2145518501952, 19.43427568696255, 23.38276245623786, 27.88654683928992, 13.39721996079594, 19.25487172391352, 18.10628386437771) },
    @{ Row = 19; Values = @(3.447932379411391, 4.214501799313131, 19.35069138706132, 23.37652773526221, 27.86988423798121, 13.40166668512386, 19.2549911578506, 18.04615509828275) },
    @{ Row = 20; Values = @(3.478292128257864, 4.214763545721846, 19.72413386078688, 23.40543176057383, 27.94618674726458, 13.38197944668573, 19.25508949297089, 18.31521720911164) },
    @{ Row = 21; Values = @(3.57805350758924, 4.216280571643483, 20.93181878274053, 23.51736004577589, 28.22547323359825, 13.32153134560464, 19.26659796512967, 19.192115042949) },
    @{ Row = 22; Values = @(3.641510528001388, 4.217745591662432, 21.68624252659342, 23.60152748223809, 28.42485744877605, 13.28631628805051, 19.28250629420882, 19.74452650928123) },
    @{ Row = 23; Values = @(3.607806124551194, 4.216920458374907, 21.28677100276507, 23.55560022602671, 28.31691661388552, 13.30471575310229, 19.27324511317518, 19.45161571405421) },
    @{ Row = 24; Values = @(3.476633165228167, 4.214746768922387, 19.70380358677297, 23.40378872247748, 27.94190984194444, 13.3830392516885, 19.2550421972668, 18.30054225635958) },
    @{ Row = 25; Values = @(3.328136223498967, 4.214445922809487, 17.84541308699434, 23.28743434430409, 27.61109985452397, 13.48565989010844, 19.27102540784917, 16.97345737691544) }
)

$cols = @("C", "D", "E", "F", "G", "H", "I", "K")

foreach ($entry in $newValues) {
    $r = $entry.Row
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "{0}{1}" -f $cols[$i], $r
        $ws.Range($addr).Value = $entry.Values[$i]
    }
}
